$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": rows 4-7 (the three "Ready for handoff" /
# "low" priority files plus the newly-promoted one) move from Priority
# "low" to "ht" (handoff-triggered) on both the zh-cn and de-de sheets,
# and their "Latest Handoff Datetime" timestamps are refreshed to the
# new handoff-generation time for each locale.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$zhCnHandoffTime = "2016-09-05 16:37:33"
$deDeHandoffTime = "2016-09-05 16:37:38"

foreach ($row in 4..7) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = $zhCnHandoffTime

    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = $deDeHandoffTime

    # Overview's "Latest HO Xliff Generate Date" column mirrors de-de's
    # Latest Handoff Datetime (they share the same underlying string).
    $wsOverview.Range("G$row").Value = $deDeHandoffTime
}
